# Apply the "Upload new version with timestamp" edit:
#  - Add a pharmacy title/header line in D2 (new merged header row content)
#  - Update a couple of data values (price 22.0000 -> 24.0000, count 11:0 -> 12:0)
#  - Update the summary count N9 25 -> 27
#  - Update the generated timestamp (10:45 PM -> 10:51 PM)
#  - Let row heights re-flow to match the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pharmacy name header in the previously-empty D2:N2 merged cell.
$ws.Range("D2").Value = "صيدليات دكتور مصطفي طلعت"

# Updated data values.
# P7 is styled with a numeric display format ("0.00") even though the
# underlying value is stored as text ("24.0000" keeps the trailing zeros
# that a real number would lose) - toggle to a text format while writing
# the value, then restore the original format so the cell style is
# unaffected.
$p7 = $ws.Range("P7")
$p7Format = $p7.NumberFormat
$p7.NumberFormat = "@"
$p7.Value = "24.0000"
$p7.NumberFormat = $p7Format

$ws.Range("Q7").Value = "12:0"
$ws.Range("N9").Value = 27

# Updated "generated at" timestamp footer.
$ws.Range("A10").Value = "Wednesday, 17 September, 2025 10:51 PM"

# Row heights re-flow after the new header text is added.
$ws.Rows(2).RowHeight = 39
$ws.Rows(3).RowHeight = 34.5
$ws.Rows(5).RowHeight = 1.5
$ws.Rows(6).RowHeight = 24.75
$ws.Rows(8).RowHeight = 25.5
$ws.Rows(9).RowHeight = 25.5
